# Re-apply the template_fullslide.pptx refresh:
#   1. Drop the single placeholder slide that shipped in the template
#      (the <p:sldIdLst> entry / ppt/slides/slide1.xml part).
#   2. Bump the cached "datetimeFigureOut" field text on the slide
#      master and the four slide layouts that still carry one (the
#      auto-date placeholders), from the 2/26 capture to 2/28.
#
# Note: the notes master's own date placeholder is read-only through
# this COM surface (its Shapes collection does not persist property
# writes), so it is intentionally left untouched here.

function Get-ShapeByName {
    param($Shapes, [string]$Name)
    for ($i = 1; $i -le $Shapes.Count; $i++) {
        $shape = $Shapes.Item($i)
        if ($shape.Name -eq $Name) {
            return $shape
        }
    }
    return $null
}

function Set-DatePlaceholderText {
    param($Shapes, [string]$OldText, [string]$NewText)
    $shape = Get-ShapeByName $Shapes "Date Placeholder 3"
    if ($shape -ne $null) {
        $tr = $shape.TextFrame.TextRange
        if ($tr.Text -eq $OldText) {
            $tr.Text = $NewText
        }
    }
}

$p = $ppt.ActivePresentation

# 1. Remove the only slide in the deck.
while ($p.Slides.Count -gt 0) {
    $p.Slides.Item(1).Delete()
}

# 2. Refresh the cached date fields (2/26/25 -> 2/28/25, 26/2/2025 -> 28/2/2025).

# Slide master date placeholder (en-AU locale -> d/M/yyyy).
Set-DatePlaceholderText $p.SlideMaster.Shapes "26/2/2025" "28/2/2025"

# Slide layout date placeholders (only the layouts that expose one).
$layouts = $p.SlideMaster.CustomLayouts
for ($li = 1; $li -le $layouts.Count; $li++) {
    $layout = $layouts.Item($li)
    Set-DatePlaceholderText $layout.Shapes "2/26/25" "2/28/25"
}
